$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.04558966666666667
$ws.Range("H2").Value = 0.136769
$ws.Range("I2").Value = 0.02375599288687187
$ws.Range("J2").Value = 0.02375599288687187
$ws.Range("M2").Value = 39.327127
$ws.Range("N2").Value = 117.981381
$ws.Range("O2").Value = 0.9923865713449503
$ws.Range("P2").Value = 0.9923865713449502
$ws.Range("Q2").Value = 1.792910610887666
$ws.Range("R2").Value = 16.136195497989
$ws.Range("S2").Value = 0.0235751283298978
$ws.Range("T2").Value = 0.0235751283298978
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.04558966666666667
$ws.Range("H3").Value = 0.136769
$ws.Range("I3").Value = 0.02375599288687187
$ws.Range("J3").Value = 0.02375599288687187
$ws.Range("O3").Value = 0.001455135597170125
$ws.Range("P3").Value = 0.001455135597170125
$ws.Range("Q3").Value = 0.002628943324888889
$ws.Range("R3").Value = 0.023660489924
$ws.Range("S3").Value = 0.00003456819089580755
$ws.Range("T3").Value = 0.00003456819089580754
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.04558966666666667
$ws.Range("H4").Value = 0.136769
$ws.Range("I4").Value = 0.02375599288687187
$ws.Range("J4").Value = 0.02375599288687187
$ws.Range("O4").Value = 0.00615829305787961
$ws.Range("P4").Value = 0.006158293057879609
$ws.Range("Q4").Value = 0.01112597579133333
$ws.Range("R4").Value = 0.100133782122
$ws.Range("S4").Value = 0.0001462963660782604
$ws.Range("T4").Value = 0.0001462963660782604
$ws.Range("I5").Value = 0.1978186777627204
$ws.Range("J5").Value = 0.1978186777627204
$ws.Range("M5").Value = 39.327127
$ws.Range("N5").Value = 117.981381
$ws.Range("O5").Value = 0.9923865713449503
$ws.Range("P5").Value = 0.9923865713449502
$ws.Range("Q5").Value = 14.92975722301
$ws.Range("R5").Value = 134.36781500709
$ws.Range("S5").Value = 0.1963125993729377
$ws.Range("T5").Value = 0.1963125993729376
$ws.Range("I6").Value = 0.1978186777627204
$ws.Range("J6").Value = 0.1978186777627204
$ws.Range("O6").Value = 0.001455135597170125
$ws.Range("P6").Value = 0.001455135597170125
$ws.Range("S6").Value = 0.0002878529997976607
$ws.Range("T6").Value = 0.0002878529997976606
$ws.Range("I7").Value = 0.1978186777627204
$ws.Range("J7").Value = 0.1978186777627204
$ws.Range("O7").Value = 0.00615829305787961
$ws.Range("P7").Value = 0.006158293057879609
$ws.Range("S7").Value = 0.001218225389985084
$ws.Range("T7").Value = 0.001218225389985084
$ws.Range("I8").Value = 0.7784253293504076
$ws.Range("J8").Value = 0.7784253293504078
$ws.Range("M8").Value = 39.327127
$ws.Range("N8").Value = 117.981381
$ws.Range("O8").Value = 0.9923865713449503
$ws.Range("P8").Value = 0.9923865713449502
$ws.Range("Q8").Value = 58.749261267347
$ws.Range("R8").Value = 528.7433514061231
$ws.Range("S8").Value = 0.7724988436421147
$ws.Range("T8").Value = 0.7724988436421149
$ws.Range("I9").Value = 0.7784253293504076
$ws.Range("J9").Value = 0.7784253293504078
$ws.Range("O9").Value = 0.001455135597170125
$ws.Range("P9").Value = 0.001455135597170125
$ws.Range("S9").Value = 0.001132714406476657
$ws.Range("T9").Value = 0.001132714406476657
$ws.Range("I10").Value = 0.7784253293504076
$ws.Range("J10").Value = 0.7784253293504078
$ws.Range("O10").Value = 0.00615829305787961
$ws.Range("P10").Value = 0.006158293057879609
$ws.Range("S10").Value = 0.004793771301816264
$ws.Range("T10").Value = 0.004793771301816264